# working on counting directors
# Populate the Title (B) and Year (D) columns for the first batch of films,
# and add two new summary rows at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$films = @(
    @{ Row = 2;  Title = "Alien";                                  Year = 1979 },
    @{ Row = 3;  Title = "Predator";                                Year = 1987 },
    @{ Row = 4;  Title = "Inception";                               Year = 2010 },
    @{ Row = 5;  Title = "Aliens";                                  Year = 1986 },
    @{ Row = 6;  Title = "The Lord of the Rings: The Two Towers";   Year = 2002 },
    @{ Row = 7;  Title = "Collateral";                              Year = 2004 },
    @{ Row = 8;  Title = "Batman Begins";                           Year = 2005 },
    @{ Row = 9;  Title = "The Dark Knight";                         Year = 2008 },
    @{ Row = 10; Title = "The Blues Brothers";                      Year = 1980 }
)

foreach ($film in $films) {
    $ws.Cells.Item($film.Row, 2).Value = $film.Title
    $ws.Cells.Item($film.Row, 4).Value = $film.Year
}

# New summary labels at the bottom of the sheet, styled like the "Year"
# header (same dark-blue fill) but with a new light-yellow font color.
$ws.Range("A103").Value = "Year with most Films"
$ws.Range("A104").Value = "Most Films by Director"

$ws.Range("D1").Copy()
$ws.Range("A103:A104").PasteSpecial(-4122)
$ws.Range("A103:A104").Font.Color = 10092543
$excel.CutCopyMode = $false

Write-Host "edit complete"
